$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Anchor cell with the default (unstyled) format, used to reset style
# after forcing a numeric-looking string to stay text via a quote prefix.
$defaultStyle = $ws.Range("D4").Style

$ws.Range('D2').Value = '27.552.05'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').Value = '1.583.14'
$ws.Range('E3').Value = '  -0.86%  '
$ws.Range('D5').Value = '''208.56'
$ws.Range('D5').Style = $defaultStyle
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('D8').Value = '''22.41'
$ws.Range('D8').Style = $defaultStyle
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('E9').Value = '  -0.95%  '
$ws.Range('E10').Value = '  -0.22%  '
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('D12').Value = '1.807.86'
$ws.Range('E12').Value = '  -0.94%  '
$ws.Range('D13').Value = '1.571.29'
$ws.Range('E13').Value = '  -1.89%  '
$ws.Range('D14').Value = '''3.84'
$ws.Range('D14').Style = $defaultStyle
$ws.Range('E14').Value = '  -0.89%  '
$ws.Range('D15').Value = '''0.526'
$ws.Range('D15').Style = $defaultStyle
$ws.Range('E15').Value = '  -2.20%  '
$ws.Range('D16').Value = '27.577.33'
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('D17').Value = '''63.13'
$ws.Range('D17').Style = $defaultStyle
$ws.Range('E17').Value = '  -0.62%  '
$ws.Range('D18').Value = '''215.51'
$ws.Range('D18').Style = $defaultStyle
$ws.Range('E18').Value = '  -0.99%  '
$ws.Range('D19').Value = '''7.33'
$ws.Range('D19').Style = $defaultStyle
$ws.Range('E19').Value = '  -0.85%  '
$ws.Range('D20').Value = '0.0₃0692'
$ws.Range('E20').Value = '  -0.49%  '
$ws.Range('E21').Value = '  -0.19%  '
$ws.Range('E22').Value = '  -1.06%  '
$ws.Range('D23').Value = '''9.80'
$ws.Range('D23').Style = $defaultStyle
$ws.Range('E23').Value = '  +1.00%  '
$ws.Range('D24').Value = '''2.01'
$ws.Range('D24').Style = $defaultStyle
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').Value = '''153.27'
$ws.Range('D25').Style = $defaultStyle
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('D26').Value = '''6.95'
$ws.Range('D26').Style = $defaultStyle
$ws.Range('E26').Value = '  +3.10%  '
$ws.Range('E27').Value = '  -0.21%  '
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('E29').Value = '  -1.39%  '
$ws.Range('E30').Value = '  -0.35%  '
$ws.Range('E31').Value = '  +0.76%  '
$ws.Range('E32').Value = '  -1.08%  '
$ws.Range('D33').Value = '1.374.94'
$ws.Range('E33').Value = '  +0.47%  '
$ws.Range('E34').Value = '  -0.19%  '
$ws.Range('D35').Value = '''1.56'
$ws.Range('D35').Style = $defaultStyle
$ws.Range('E35').Value = '  +1.43%  '
$ws.Range('D36').Value = '''0.974'
$ws.Range('D36').Style = $defaultStyle
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('E38').Value = '  +1.37%  '
$ws.Range('E39').Value = '  -1.17%  '
$ws.Range('D40').Value = '''0.829'
$ws.Range('D40').Style = $defaultStyle
$ws.Range('E40').Value = '  +1.83%  '
$ws.Range('E41').Value = '  -0.22%  '
$ws.Range('D42').Value = '''0.971'
$ws.Range('D42').Style = $defaultStyle
$ws.Range('E42').Value = '  -0.48%  '
$ws.Range('E43').Value = '  +1.03%  '
$ws.Range('D44').Value = '''64.53'
$ws.Range('D44').Style = $defaultStyle
$ws.Range('E44').Value = '  +0.64%  '
$ws.Range('D45').Value = '''5.31'
$ws.Range('D45').Style = $defaultStyle
$ws.Range('E45').Value = '  -1.19%  '
$ws.Range('E46').Value = '  +2.99%  '
$ws.Range('D47').Value = '1.721.22'
$ws.Range('E47').Value = '  -0.82%  '
$ws.Range('D48').Value = '''85.36'
$ws.Range('D48').Style = $defaultStyle
$ws.Range('E48').Value = '  -2.92%  '
$ws.Range('D49').Value = '0.0₇0995'
$ws.Range('E49').Value = '  -0.76%  '
$ws.Range('E50').Value = '  -1.26%  '
$ws.Range('D51').Value = '''0.0496'
$ws.Range('D51').Style = $defaultStyle
$ws.Range('E51').Value = '  -0.57%  '
